# New RCC scripts added to the "Test Cases" sheet (4 new rows: RCC002, RCC007,
# RCC008, RCC012), mirroring the existing RCC111 row (same Description/Runmode).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Clone formatting (borders/fill/etc.) from the last existing data row block
# down into the 4 new rows before filling in values.
$ws.Range("A2:E5").Copy()
$ws.Range("A6:E9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill column B ("Jira id") first, then column A ("TCID"), matching the order
# new test-case scripts were entered.
$ws.Range("B6").Value = "ABCD1"
$ws.Range("B7").Value = "ABCD2"
$ws.Range("B8").Value = "ABCD3"
$ws.Range("B9").Value = "ABCD4"

$ws.Range("A6").Value = "RCC002"
$ws.Range("A7").Value = "RCC007"
$ws.Range("A8").Value = "RCC008"
$ws.Range("A9").Value = "RCC012"

$desc = "Verify that user is able to add an article to the group from search results page."
$ws.Range("C6").Value = $desc
$ws.Range("C7").Value = $desc
$ws.Range("C8").Value = $desc
$ws.Range("C9").Value = $desc

$ws.Range("D6").Value = "Y"
$ws.Range("D7").Value = "Y"
$ws.Range("D8").Value = "Y"
$ws.Range("D9").Value = "Y"

$ws.Range("A9").Select() | Out-Null
